$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header O2: "Utility" -> "Utility (Percent)"
$ws.Range("O2").Value = "Utility (Percent)"

# Append " msec" to I, J, K columns for data rows 3 through 23
for ($row = 3; $row -le 23; $row++) {
    $cellI = $ws.Cells.Item($row, 9)
    $cellJ = $ws.Cells.Item($row, 10)
    $cellK = $ws.Cells.Item($row, 11)
    $cellI.Value = [string]$cellI.Value2 + " msec"
    $cellJ.Value = [string]$cellJ.Value2 + " msec"
    $cellK.Value = [string]$cellK.Value2 + " msec"
}
